# "reporte generando hojas dinamicamente"
# Roll the Summary report forward by one year:
#   - Table 1 (rows 12-23) currently labelled with 2017 month-end dates -> relabel 2018
#   - Table 2 (rows 32-43) currently labelled with 2018 month-end dates -> relabel 2019
#   - Zero out the (now stale) revenue figures that belonged to the just-retired period
#   - Bump the big "2018" year banner (A26) to "2019"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextDate($addr, $text) {
    # Writing a date-shaped string straight into a date-formatted cell makes
    # Excel "helpfully" coerce it into a real date serial. Flip the cell to
    # Text, assign the literal string, then restore its original display
    # format so neither the style nor the number format actually changes.
    $rng = $ws.Range($addr)
    $origFormat = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.NumberFormat = $origFormat
}

# --- Table 1 (rows 12-23): 2017 -> 2018 month-end labels ---
Set-TextDate "A12" "1/31/2018"
Set-TextDate "A13" "2/28/2018"
Set-TextDate "A14" "3/31/2018"
Set-TextDate "A15" "4/30/2018"
Set-TextDate "A16" "5/31/2018"
Set-TextDate "A17" "6/30/2018"
Set-TextDate "A18" "7/31/2018"
Set-TextDate "A19" "8/30/2018"
Set-TextDate "A20" "9/31/2018"
Set-TextDate "A21" "10/30/2018"
Set-TextDate "A22" "11/31/2018"
Set-TextDate "A23" "12/30/2018"

# --- Table 2 (rows 32-43): 2018 -> 2019 month-end labels ---
Set-TextDate "A32" "1/31/2019"
Set-TextDate "A33" "2/28/2019"
Set-TextDate "A34" "3/31/2019"
Set-TextDate "A35" "4/30/2019"
Set-TextDate "A36" "5/31/2019"
Set-TextDate "A37" "6/30/2019"
Set-TextDate "A38" "7/31/2019"
Set-TextDate "A39" "8/30/2019"
Set-TextDate "A40" "9/31/2019"
Set-TextDate "A41" "10/30/2019"
Set-TextDate "A42" "11/31/2019"
Set-TextDate "A43" "12/30/2019"

# --- Year banner above table 2 ---
$ws.Range("A26").Value = 2019

# --- Zero out the retired 2017 figures in table 1 ---
$ws.Range("B12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("K12").Value = 1
$ws.Range("R12").Value = 0

$ws.Range("B13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("R13").Value = 0

$ws.Range("B14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("R14").Value = 0

$ws.Range("B15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("R15").Value = 0

$ws.Range("B16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("R16").Value = 0

$ws.Range("B17").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("R17").Value = 0

$ws.Range("R18").Value = 0
$ws.Range("R19").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("R23").Value = 0
